# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, per the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> list of (row, newValue)
$updates = @{
    "展览"     = @(
        @{ Row = 2;  Value = 1656 },
        @{ Row = 4;  Value = 9485 },
        @{ Row = 11; Value = 1640 },
        @{ Row = 12; Value = 1415 },
        @{ Row = 15; Value = 1474 },
        @{ Row = 21; Value = 387 },
        @{ Row = 35; Value = 181 },
        @{ Row = 38; Value = 254 },
        @{ Row = 43; Value = 318 }
    )
    "演出"     = @(
        @{ Row = 11; Value = 693 },
        @{ Row = 22; Value = 284 },
        @{ Row = 36; Value = 115 }
    )
    "本地生活" = @(
        @{ Row = 7;  Value = 2380 },
        @{ Row = 8;  Value = 3609 },
        @{ Row = 11; Value = 73 }
    )
    "全部类型" = @(
        @{ Row = 2;  Value = 1656 },
        @{ Row = 3;  Value = 9485 },
        @{ Row = 6;  Value = 3609 },
        @{ Row = 8;  Value = 73 },
        @{ Row = 9;  Value = 73 },
        @{ Row = 13; Value = 693 },
        @{ Row = 14; Value = 1415 },
        @{ Row = 18; Value = 1474 },
        @{ Row = 33; Value = 284 },
        @{ Row = 43; Value = 254 },
        @{ Row = 48; Value = 318 },
        @{ Row = 49; Value = 115 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
